$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2693.2144
$ws.Range("I138").Value = 2558.8333
$ws.Range("K138").Value = 7676.499899999999
$ws.Range("M138").Value = -2536.499899999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 50550
$ws.Range("J92").Value = 50550
$ws.Range("L92").Value = 50550
$ws.Range("N92").Value = -55542
$ws.Range("H97").Value = 620.2917
$ws.Range("I97").Value = 554.1429000000001
$ws.Range("J97").Value = 1083.3334
$ws.Range("K97").Value = 554.1429000000001
$ws.Range("L97").Value = 1083.3334
$ws.Range("M97").Value = -58.14290000000005
$ws.Range("N97").Value = -2075.3334

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2999.5652
$ws.Range("I86").Value = 2404.6428
$ws.Range("J86").Value = 3925
$ws.Range("K86").Value = 2404.6428
$ws.Range("L86").Value = 3925
$ws.Range("M86").Value = -1281.6428
$ws.Range("N86").Value = -6171
$ws.Range("H89").Value = 2999.5652
$ws.Range("I89").Value = 2404.6428
$ws.Range("J89").Value = 3925
$ws.Range("K89").Value = 12023.214
$ws.Range("L89").Value = 19625
$ws.Range("M89").Value = -6407.214
$ws.Range("N89").Value = -30857
$ws.Range("H94").Value = 2554.8333
$ws.Range("I94").Value = 1991.5
$ws.Range("J94").Value = 4244.8335
$ws.Range("K94").Value = 1991.5
$ws.Range("L94").Value = 4244.8335
$ws.Range("M94").Value = -1540.5
$ws.Range("N94").Value = -5146.8335
$ws.Range("H105").Value = 5750
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").Value = ""
$ws.Range("H134").Value = 1379.8889
$ws.Range("I134").Value = 1284.5883
$ws.Range("K134").Value = 3853.7649
$ws.Range("M134").Value = -1318.7649

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = ""
$ws.Range("N16").Value = ""
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").Value = ""
$ws.Range("H58").Value = 1708.8889
$ws.Range("I58").Value = 1549.375
$ws.Range("K58").Value = 1549.375
$ws.Range("M58").Value = -1346.375
$ws.Range("H62").Value = 2750
$ws.Range("I62").Value = 2750
$ws.Range("K62").Value = 2750
$ws.Range("M62").Value = -2126
$ws.Range("H65").Value = 2750
$ws.Range("I65").Value = 2750
$ws.Range("K65").Value = 13750
$ws.Range("M65").Value = -10630
$ws.Range("H105").Value = 2304.125
$ws.Range("I105").Value = 2186.8
$ws.Range("K105").Value = 2186.8
$ws.Range("M105").Value = -439.8000000000002
$ws.Range("H107").Value = 308
$ws.Range("I107").Value = 315
$ws.Range("J107").Value = 280
$ws.Range("K107").Value = 315
$ws.Range("L107").Value = 280
$ws.Range("M107").Value = 1605
$ws.Range("N107").Value = -4120
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = ""
$ws.Range("H134").Value = 2890.2856
$ws.Range("I134").Value = 2505.818
$ws.Range("J134").Value = 4300
$ws.Range("K134").Value = 7517.454000000001
$ws.Range("L134").Value = 12900
$ws.Range("M134").Value = -4982.454000000001
$ws.Range("N134").Value = -17970
$ws.Range("H136").Value = 1708.8889
$ws.Range("I136").Value = 1549.375
$ws.Range("K136").Value = 4648.125
$ws.Range("M136").Value = -2098.125

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 11054
$ws.Range("I80").Value = 2499
$ws.Range("J80").Value = 13498.286
$ws.Range("K80").Value = 7497
$ws.Range("L80").Value = 40494.858
$ws.Range("M80").Value = -6561
$ws.Range("N80").Value = -42366.858
$ws.Range("H83").Value = 11054
$ws.Range("I83").Value = 2499
$ws.Range("J83").Value = 13498.286
$ws.Range("K83").Value = 22491
$ws.Range("L83").Value = 121484.574
$ws.Range("M83").Value = -17811
$ws.Range("N83").Value = -130844.574
$ws.Range("H92").Value = 637.2222
$ws.Range("I92").Value = 375.5
$ws.Range("K92").Value = 1126.5
$ws.Range("M92").Value = 121.5
$ws.Range("H120").Value = 5412.857
$ws.Range("I120").Value = 2472.5
$ws.Range("J120").Value = 9333.333000000001
$ws.Range("K120").Value = 7417.5
$ws.Range("L120").Value = 27999.999
$ws.Range("M120").Value = -2579.5
$ws.Range("N120").Value = -37675.999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2000
$ws.Range("I7").Value = 2000
$ws.Range("K7").Value = 2000
$ws.Range("M7").Value = -1888
$ws.Range("H16").Value = 497.22223
$ws.Range("I16").Value = 497.22223
$ws.Range("K16").Value = 497.22223
$ws.Range("M16").Value = -327.22223
$ws.Range("H46").Value = 1629.5454
$ws.Range("I46").Value = 1241.6666
$ws.Range("J46").Value = 2308.3333
$ws.Range("K46").Value = 1241.6666
$ws.Range("L46").Value = 2308.3333
$ws.Range("M46").Value = -1053.6666
$ws.Range("N46").Value = -2684.3333
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530
$ws.Range("H132").Value = 3044.1052
$ws.Range("I132").Value = 2989.2666
$ws.Range("J132").Value = 3249.75
$ws.Range("K132").Value = 8967.799800000001
$ws.Range("L132").Value = 9749.25
$ws.Range("M132").Value = -6437.799800000001
$ws.Range("N132").Value = -14809.25
$ws.Range("H136").Value = 3867.8572
$ws.Range("I136").Value = 3015.1
$ws.Range("J136").Value = 5999.75
$ws.Range("K136").Value = 9045.299999999999
$ws.Range("L136").Value = 17999.25
$ws.Range("M136").Value = -6495.299999999999
$ws.Range("N136").Value = -23099.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2044.1111
$ws.Range("I96").Value = 1399.5
$ws.Range("K96").Value = 1399.5
$ws.Range("M96").Value = -26.5
$ws.Range("H113").Value = 450.1111
$ws.Range("J113").Value = 683
$ws.Range("L113").Value = 2049
$ws.Range("N113").Value = -6389
